$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 2742.8572
$ws.Range("I76").Value = 2742.8572
$ws.Range("K76").Value = 2742.8572
$ws.Range("M76").Value = -2427.8572
$ws.Range("H79").Value = 2742.8572
$ws.Range("I79").Value = 2742.8572
$ws.Range("K79").Value = 2742.8572
$ws.Range("M79").Value = -1650.8572
$ws.Range("H141").Value = 14459
$ws.Range("I141").Value = 18023.625
$ws.Range("J141").Value = 4953.3335
$ws.Range("K141").Value = 54070.875
$ws.Range("L141").Value = 14860.0005
$ws.Range("M141").Value = -48890.875
$ws.Range("N141").Value = -25220.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1966.7778
$ws.Range("I88").Value = 1528.7142
$ws.Range("J88").Value = 3500
$ws.Range("K88").Value = 1528.7142
$ws.Range("L88").Value = 3500
$ws.Range("M88").Value = -1122.7142
$ws.Range("N88").Value = -4312
$ws.Range("H91").Value = 1966.7778
$ws.Range("I91").Value = 1528.7142
$ws.Range("J91").Value = 3500
$ws.Range("K91").Value = 1528.7142
$ws.Range("L91").Value = 3500
$ws.Range("M91").Value = -124.7141999999999
$ws.Range("N91").Value = -6308

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H82").Value = 33141.5
$ws.Range("I82").Value = 20000
$ws.Range("J82").Value = 46283
$ws.Range("K82").Value = 20000
$ws.Range("L82").Value = 46283
$ws.Range("M82").Value = -19617
$ws.Range("N82").Value = -47049
$ws.Range("H85").Value = 33141.5
$ws.Range("I85").Value = 20000
$ws.Range("J85").Value = 46283
$ws.Range("K85").Value = 20000
$ws.Range("L85").Value = 46283
$ws.Range("M85").Value = -18674
$ws.Range("N85").Value = -48935
$ws.Range("H86").Value = 3015.24
$ws.Range("I86").Value = 3081.1462
$ws.Range("J86").Value = 2715
$ws.Range("K86").Value = 3081.1462
$ws.Range("L86").Value = 2715
$ws.Range("M86").Value = -1958.1462
$ws.Range("N86").Value = -4961
$ws.Range("H89").Value = 3015.24
$ws.Range("I89").Value = 3081.1462
$ws.Range("J89").Value = 2715
$ws.Range("K89").Value = 15405.731
$ws.Range("L89").Value = 13575
$ws.Range("M89").Value = -9789.731
$ws.Range("N89").Value = -24807

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1242.5807
$ws.Range("I16").Value = 1179.6957
$ws.Range("J16").Value = 1423.375
$ws.Range("K16").Value = 1179.6957
$ws.Range("L16").Value = 1423.375
$ws.Range("M16").Value = -892.6957
$ws.Range("N16").Value = -1997.375
$ws.Range("H62").Value = 3604.25
$ws.Range("I62").Value = 2549.4443
$ws.Range("J62").Value = 6768.6665
$ws.Range("K62").Value = 2549.4443
$ws.Range("L62").Value = 6768.6665
$ws.Range("M62").Value = -1925.4443
$ws.Range("N62").Value = -8016.6665
$ws.Range("H65").Value = 3604.25
$ws.Range("I65").Value = 2549.4443
$ws.Range("J65").Value = 6768.6665
$ws.Range("K65").Value = 12747.2215
$ws.Range("L65").Value = 33843.3325
$ws.Range("M65").Value = -9627.2215
$ws.Range("N65").Value = -40083.3325
$ws.Range("H99").Value = 1978.5
$ws.Range("I99").Value = 1800
$ws.Range("J99").Value = 2276
$ws.Range("K99").Value = 1800
$ws.Range("L99").Value = 2276
$ws.Range("M99").Value = -302
$ws.Range("N99").Value = -5272
$ws.Range("H113").Value = 1242.5807
$ws.Range("I113").Value = 1179.6957
$ws.Range("J113").Value = 1423.375
$ws.Range("K113").Value = 1179.6957
$ws.Range("L113").Value = 1423.375
$ws.Range("M113").Value = 990.3043
$ws.Range("N113").Value = -5763.375
$ws.Range("H126").Value = 1978.5
$ws.Range("I126").Value = 1800
$ws.Range("J126").Value = 2276
$ws.Range("K126").Value = 5400
$ws.Range("L126").Value = 6828
$ws.Range("M126").Value = -2930
$ws.Range("N126").Value = -11768

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 1463.9412
$ws.Range("I60").Value = 515.8333
$ws.Range("J60").Value = 1981.091
$ws.Range("K60").Value = 1547.4999
$ws.Range("L60").Value = 5943.272999999999
$ws.Range("M60").Value = -1296.4999
$ws.Range("N60").Value = -6445.272999999999
$ws.Range("H75").Value = 2952.2727
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 2952.2727
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 8856.8181
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -10852.8181
$ws.Range("H78").Value = 2952.2727
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 2952.2727
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 26570.4543
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -36554.4543
$ws.Range("H103").Value = 2505.8125
$ws.Range("I103").Value = 1136.5625
$ws.Range("J103").Value = 3875.0625
$ws.Range("K103").Value = 3409.6875
$ws.Range("L103").Value = 11625.1875
$ws.Range("M103").Value = -2530.6875
$ws.Range("N103").Value = -13383.1875
$ws.Range("H136").Value = 7000
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 7000
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 21000
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -31200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 28623.715
$ws.Range("H80").Value = 2955
$ws.Range("I80").Value = 2875.625
$ws.Range("J80").Value = 3045.7144
$ws.Range("K80").Value = 2875.625
$ws.Range("L80").Value = 3045.7144
$ws.Range("M80").Value = -1877.625
$ws.Range("N80").Value = -5041.7144
$ws.Range("H83").Value = 2955
$ws.Range("I83").Value = 2875.625
$ws.Range("J83").Value = 3045.7144
$ws.Range("K83").Value = 14378.125
$ws.Range("L83").Value = 15228.572
$ws.Range("M83").Value = -9386.125
$ws.Range("N83").Value = -25212.572

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1145.1818
$ws.Range("I16").Value = 978.1429000000001
$ws.Range("J16").Value = 1437.5
$ws.Range("K16").Value = 978.1429000000001
$ws.Range("L16").Value = 1437.5
$ws.Range("M16").Value = -808.1429000000001
$ws.Range("N16").Value = -1777.5
$ws.Range("H55").Value = 405.41666
$ws.Range("I55").Value = 342.4
$ws.Range("J55").Value = 484.1875
$ws.Range("K55").Value = 342.4
$ws.Range("L55").Value = 484.1875
$ws.Range("M55").Value = -169.4
$ws.Range("N55").Value = -830.1875
$ws.Range("H93").Value = 2640.7693
$ws.Range("I93").Value = 2441.7144
$ws.Range("J93").Value = 3476.8
$ws.Range("K93").Value = 2441.7144
$ws.Range("L93").Value = 3476.8
$ws.Range("M93").Value = -1193.7144
$ws.Range("N93").Value = -5972.8

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 900
$ws.Range("I113").Value = 750
$ws.Range("J113").Value = 1000
$ws.Range("K113").Value = 2250
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = -80
$ws.Range("N113").Value = -7340
